$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3340
$ws.Range("E2").Value = 243
$ws.Range("F2").Value = 243
$ws.Range("G2").Value = 220
$ws.Range("H2").Value = 158
$ws.Range("I2").Value = 166
$ws.Range("J2").Value = -8
$ws.Range("K2").Value = 3257
$ws.Range("L2").Value = 2452
$ws.Range("M2").Value = 805
$ws.Range("N2").Value = 808
$ws.Range("O2").Value = -3
$ws.Range("P2").Value = 45
$ws.Range("Q2").Value = 12
$ws.Range("R2").Value = -456
$ws.Range("S2").Value = 435
$ws.Range("T2").Value = 415
$ws.Range("U2").Value = -404
$ws.Range("V2").Value = 1450
$ws.Range("W2").Value = 7.29
$ws.Range("X2").Value = 4.72
$ws.Range("AA2").Value = 304.43
$ws.Range("AB2").Value = 1678.1
$ws.Range("AC2").Value = 2170
$ws.Range("AD2").Value = 44.8
$ws.Range("AE2").Value = 8738
$ws.Range("AF2").Value = 11.13
$ws.Range("AG2").Value = 487
$ws.Range("AH2").Value = 0.5
$ws.Range("AI2").Value = 27.05
$ws.Range("AJ2").Value = 9248916
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()

# Row 3
$ws.Range("D3").Value = 5333
$ws.Range("E3").Value = 359
$ws.Range("F3").Value = 359
$ws.Range("G3").Value = 295
$ws.Range("H3").Value = 189
$ws.Range("I3").Value = 212
$ws.Range("J3").Value = -24
$ws.Range("K3").Value = 4389
$ws.Range("L3").Value = 3424
$ws.Range("M3").Value = 965
$ws.Range("N3").Value = 979
$ws.Range("O3").Value = -13
$ws.Range("P3").Value = 45
$ws.Range("Q3").Value = -65
$ws.Range("R3").Value = -471
$ws.Range("S3").Value = 602
$ws.Range("T3").Value = 435
$ws.Range("U3").Value = -500
$ws.Range("V3").Value = 2106
$ws.Range("W3").Value = 6.74
$ws.Range("X3").Value = 3.54
$ws.Range("Y3").Value = 23.78
$ws.Range("Z3").Value = 4.94
$ws.Range("AA3").Value = 354.76
$ws.Range("AB3").Value = 2039.1
$ws.Range("AC3").Value = 2297
$ws.Range("AD3").Value = 78.11
$ws.Range("AE3").Value = 10585
$ws.Range("AF3").Value = 16.95
$ws.Range("AG3").Value = 681
$ws.Range("AH3").Value = 0.38
$ws.Range("AI3").Value = 29.64
$ws.Range("AJ3").Value = 9248916

# Row 4
$ws.Range("D4").Value = 7570
$ws.Range("E4").Value = 526
$ws.Range("F4").Value = 526
$ws.Range("G4").Value = 465
$ws.Range("H4").Value = 314
$ws.Range("I4").Value = 348
$ws.Range("J4").Value = -33
$ws.Range("K4").Value = 6659
$ws.Range("L4").Value = 4586
$ws.Range("M4").Value = 2073
$ws.Range("N4").Value = 2109
$ws.Range("O4").Value = -36
$ws.Range("P4").Value = 50
$ws.Range("Q4").Value = 235
$ws.Range("R4").Value = -975
$ws.Range("S4").Value = 1249
$ws.Range("T4").Value = 896
$ws.Range("U4").Value = -661
$ws.Range("V4").Value = 2497
$ws.Range("W4").Value = 6.95
$ws.Range("X4").Value = 4.15
$ws.Range("Y4").Value = 22.53
$ws.Range("Z4").Value = 5.69
$ws.Range("AA4").Value = 221.21
$ws.Range("AB4").Value = 4162.14
$ws.Range("AC4").Value = 3752
$ws.Range("AD4").Value = 31.85
$ws.Range("AE4").Value = 20997
$ws.Range("AF4").Value = 5.69
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 0.84
$ws.Range("AI4").Value = 28.88
$ws.Range("AJ4").Value = 10049509

# Row 5
$ws.Range("D5").Value = 8840
$ws.Range("E5").Value = 351
$ws.Range("F5").Value = 351
$ws.Range("G5").Value = 253
$ws.Range("H5").Value = 155
$ws.Range("I5").Value = 189
$ws.Range("J5").Value = -33
$ws.Range("K5").Value = 9340
$ws.Range("L5").Value = 7056
$ws.Range("M5").Value = 2284
$ws.Range("N5").Value = 2181
$ws.Range("O5").Value = 103
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = -66
$ws.Range("R5").Value = -1727
$ws.Range("S5").Value = 1501
$ws.Range("T5").Value = 812
$ws.Range("U5").Value = -878
$ws.Range("V5").Value = 4052
$ws.Range("W5").Value = 3.98
$ws.Range("X5").Value = 1.76
$ws.Range("Y5").Value = 8.789999999999999
$ws.Range("Z5").Value = 1.94
$ws.Range("AA5").Value = 308.96
$ws.Range("AB5").Value = 4303.58
$ws.Range("AC5").Value = 1876
$ws.Range("AD5").Value = 62.36
$ws.Range("AE5").Value = 21712
$ws.Range("AF5").Value = 5.39
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 0.26
$ws.Range("AI5").Value = 15.98
$ws.Range("AJ5").Value = 10049509

# Row 6
$ws.Range("D6").Value = 12597
$ws.Range("E6").Value = 523
$ws.Range("F6").Value = 523
$ws.Range("G6").Value = 362
$ws.Range("H6").Value = 211
$ws.Range("I6").Value = 326
$ws.Range("K6").Value = 10615
$ws.Range("L6").Value = 8197
$ws.Range("M6").Value = 2418
$ws.Range("N6").Value = 2654
$ws.Range("P6").Value = 50
$ws.Range("Q6").Value = -546
$ws.Range("R6").Value = -567
$ws.Range("S6").Value = 1075
$ws.Range("T6").Value = 490
$ws.Range("U6").Value = -1036
$ws.Range("V6").Value = 5165
$ws.Range("W6").Value = 4.15
$ws.Range("X6").Value = 1.67
$ws.Range("Y6").Value = 13.49
$ws.Range("Z6").Value = 2.11
$ws.Range("AA6").Value = 338.93
$ws.Range("AB6").Value = 5247.6
$ws.Range("AC6").Value = 3246
$ws.Range("AD6").Value = 40.05
$ws.Range("AE6").Value = 26420
$ws.Range("AF6").Value = 4.92
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 0.46
$ws.Range("AI6").Value = 18.48
$ws.Range("AJ6").Value = 10049509

# Row 7
$ws.Range("D7").Value = 13193
$ws.Range("E7").Value = 483
$ws.Range("G7").Value = 358
$ws.Range("H7").Value = 199
$ws.Range("I7").Value = 289
$ws.Range("K7").Value = 11284
$ws.Range("L7").Value = 8620
$ws.Range("M7").Value = 2664
$ws.Range("N7").Value = 2883
$ws.Range("P7").Value = 50
$ws.Range("Q7").Value = 700
$ws.Range("R7").Value = -481
$ws.Range("S7").Value = -1
$ws.Range("T7").Value = 483
$ws.Range("U7").Value = 247
$ws.Range("W7").Value = 3.66
$ws.Range("X7").Value = 1.51
$ws.Range("Y7").Value = 10.43
$ws.Range("Z7").Value = 1.82
$ws.Range("AA7").Value = 323.54
$ws.Range("AC7").Value = 2873
$ws.Range("AD7").Value = 23.8
$ws.Range("AE7").Value = 28700
$ws.Range("AF7").Value = 2.38
$ws.Range("AG7").Value = 656
$ws.Range("AH7").Value = 0.96
$ws.Range("AI7").Value = 22.84

# Row 8
$ws.Range("D8").Value = 14609
$ws.Range("E8").Value = 664
$ws.Range("G8").Value = 520
$ws.Range("H8").Value = 330
$ws.Range("I8").Value = 435
$ws.Range("K8").Value = 12153
$ws.Range("L8").Value = 9193
$ws.Range("M8").Value = 2960
$ws.Range("N8").Value = 3243
$ws.Range("P8").Value = 50
$ws.Range("Q8").Value = 460
$ws.Range("R8").Value = -452
$ws.Range("S8").Value = 136
$ws.Range("T8").Value = 363
$ws.Range("U8").Value = 40
$ws.Range("W8").Value = 4.55
$ws.Range("X8").Value = 2.26
$ws.Range("Y8").Value = 14.2
$ws.Range("Z8").Value = 2.82
$ws.Range("AA8").Value = 310.56
$ws.Range("AC8").Value = 4328
$ws.Range("AD8").Value = 15.8
$ws.Range("AE8").Value = 32286
$ws.Range("AF8").Value = 2.12
$ws.Range("AG8").Value = 778
$ws.Range("AH8").Value = 1.14
$ws.Range("AI8").Value = 17.98

# Row 9
$ws.Range("D9").Value = 16113
$ws.Range("E9").Value = 846
$ws.Range("G9").Value = 698
$ws.Range("H9").Value = 438
$ws.Range("I9").Value = 557
$ws.Range("K9").Value = 12978
$ws.Range("L9").Value = 9613
$ws.Range("M9").Value = 3365
$ws.Range("N9").Value = 3727
$ws.Range("P9").Value = 50
$ws.Range("Q9").Value = 585
$ws.Range("R9").Value = -380
$ws.Range("S9").Value = 33
$ws.Range("T9").Value = 340
$ws.Range("U9").Value = 210
$ws.Range("W9").Value = 5.25
$ws.Range("X9").Value = 2.72
$ws.Range("Y9").Value = 15.99
$ws.Range("Z9").Value = 1.82
$ws.Range("AA9").Value = 285.66
$ws.Range("AC9").Value = 5546
$ws.Range("AD9").Value = 12.33
$ws.Range("AE9").Value = 37098
$ws.Range("AF9").Value = 1.84
$ws.Range("AG9").Value = 803
$ws.Range("AH9").Value = 1.17
$ws.Range("AI9").Value = 14.48
